$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from "Checklist" to "Session"
$ws.Name = "Session"

# Row 2: log entry type changes from "Selection" to "Scan"
$ws.Range("E2").Value = "Scan"

# The 201177 log entry (old row 3) is removed entirely; everything below
# (old row 4, the 201180 entry) shifts up to become the new row 3. Deleting
# this way (rather than overwriting values) preserves the original cell's
# text formatting for the numeric-looking student ID.
$ws.Rows.Item(3).Delete()

# New row 3 (previously row 4): log time + type are updated
$ws.Range("D3").Value = "13:20:07"
$ws.Range("E3").Value = "Scan"
